# Retrocedi ações e aprimorei o index para notebook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant values shared by every data row
$operadora = "422380 - YOU ASSISTÊNCIA MÉDICA LTDA"
$hoje      = "13-04-2023"
$prazo     = "10 dias úteis"
$respondido = "NO"
$natureza   = "Assistencial"
$opcoes     = "Responder  Detalhes"

# Per-row data: DataNotificacao, Demanda, Protocolo, Beneficiario
$rows = @(
  @{ Row=2;  Data="12/04/2023  11:59:41"; Demanda=12166383; Protocolo=8604036; Benef="MARLENE NUNES HONDA TAVARES" },
  @{ Row=3;  Data="12/04/2023  12:03:08"; Demanda=12166396; Protocolo=8604069; Benef="MARCO ANTONIO DALPRA" },
  @{ Row=4;  Data="12/04/2023  15:43:18"; Demanda=12167333; Protocolo=8605179; Benef="RAVI SCHULZ XAVIER DA CRUZ" },
  @{ Row=5;  Data="12/04/2023  16:14:52"; Demanda=12167445; Protocolo=8605345; Benef="KAYKY BRUNNO SOUZA LOPES" },
  @{ Row=6;  Data="13/04/2023  08:23:06"; Demanda=12168095; Protocolo=8606114; Benef="MÔNICA ALVES GOMES" },
  @{ Row=7;  Data="13/04/2023  11:42:27"; Demanda=12168687; Protocolo=8606848; Benef="MARCIO CANDIDO DE OLIVEIRA" },
  @{ Row=8;  Data="13/04/2023  13:26:27"; Demanda=12169086; Protocolo=8607367; Benef="BARBARA KELLY CARNEIRO LEÃO RODRIGUES" },
  @{ Row=9;  Data="13/04/2023  15:12:39"; Demanda=12169516; Protocolo=8607900; Benef="ANALIS SOARES SILVA" },
  @{ Row=10; Data="13/04/2023  16:38:38"; Demanda=12169806; Protocolo=8608265; Benef="MAICKSON CAIQUE VENANCIO" },
  @{ Row=11; Data="13/04/2023  16:46:57"; Demanda=12169830; Protocolo=8608288; Benef="EMILLE FERNANDES CORREA" },
  @{ Row=12; Data="13/04/2023  16:58:57"; Demanda=12169862; Protocolo=8608371; Benef="MILENA FREIRE TRAVASSOS COUSSEIRO" },
  @{ Row=13; Data="13/04/2023  18:17:42"; Demanda=12170051; Protocolo=8608573; Benef="ROGERIA DORALICE SOARES DA SILVA" }
)

foreach ($r in $rows) {
  $row = $r.Row
  $ws.Cells.Item($row, 1).Value = $hoje
  $ws.Cells.Item($row, 2).Value = $operadora
  $ws.Cells.Item($row, 3).Value = $r.Data
  $ws.Cells.Item($row, 4).Value = $r.Demanda
  $ws.Cells.Item($row, 5).Value = $r.Protocolo
  $ws.Cells.Item($row, 6).Value = $r.Benef
  $ws.Cells.Item($row, 7).Value = $prazo
  $ws.Cells.Item($row, 8).Value = $respondido
  $ws.Cells.Item($row, 9).Value = $natureza
  $ws.Cells.Item($row, 10).Value = $opcoes
}
